$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.05003266666666667
$ws.Range("H2").Value = 0.150098
$ws.Range("M2").Value = 5.616015666666666
$ws.Range("N2").Value = 16.848047
$ws.Range("O2").Value = 0.2860808099623356
$ws.Range("P2").Value = 0.2860808099623357
$ws.Range("Q2").Value = 0.2809842398451111
$ws.Range("R2").Value = 2.528858158606
$ws.Range("S2").Value = 0.2860808099623356
$ws.Range("T2").Value = 0.2860808099623357

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.05003266666666667
$ws.Range("H3").Value = 0.150098
$ws.Range("M3").Value = 8.435525999999999
$ws.Range("N3").Value = 25.306578
$ws.Range("O3").Value = 0.4297071542841152
$ws.Range("P3").Value = 0.4297071542841153
$ws.Range("Q3").Value = 0.422051860516
$ws.Range("R3").Value = 3.798466744644
$ws.Range("S3").Value = 0.4297071542841152
$ws.Range("T3").Value = 0.4297071542841153

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.05003266666666667
$ws.Range("H4").Value = 0.150098
$ws.Range("M4").Value = 2.036951
$ws.Range("N4").Value = 6.110853000000001
$ws.Range("O4").Value = 0.1037626364528048
$ws.Range("P4").Value = 0.1037626364528048
$ws.Range("Q4").Value = 0.1019140903993333
$ws.Range("R4").Value = 0.9172268135940002
$ws.Range("S4").Value = 0.1037626364528048
$ws.Range("T4").Value = 0.1037626364528048

# Row 5
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.05003266666666667
$ws.Range("H5").Value = 0.150098
$ws.Range("M5").Value = 3.542379
$ws.Range("N5").Value = 10.627137
$ws.Range("O5").Value = 0.1804493993007442
$ws.Range("P5").Value = 0.1804493993007443
$ws.Range("Q5").Value = 0.177234667714
$ws.Range("R5").Value = 1.595112009426
$ws.Range("S5").Value = 0.1804493993007442
$ws.Range("T5").Value = 0.1804493993007443
